# Closes #5044: replace the "single underlines for emphasis" sample
# paragraph with a plain "underlining" example (bare <w:u/> instead of
# <w:u w:val="single"/>, and the stray <w:i/> emphasis run removed).

$d = $word.ActiveDocument

# Find the paragraph that starts with the old sample text. Matching on
# text (rather than a hard-coded paragraph index) keeps this robust to
# unrelated structural changes elsewhere in the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Some people use single underlines for emphasis.")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found; no changes made."
} else {
    $r = $target.Range

    # Preserve the paragraph's own attributes (w14:paraId / rsids / etc.)
    # by reading them back off the round-tripped OOXML instead of
    # hard-coding them.
    $openxml = $r.WordOpenXML
    $pTagAttrs = ""
    if ($openxml -match "<w:p([^>]*)>") {
        $pTagAttrs = $matches[1]
    }

    # Five runs: "Some people use" | " " | "underlining" (bare <w:u/>) |
    # " " | "for emphasis." -- no italics anywhere any more.
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"' + $pTagAttrs + '>' +
        '<w:r><w:t xml:space="preserve">Some people use</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:u/></w:rPr><w:t xml:space="preserve">underlining</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">for emphasis.</w:t></w:r>' +
        '</w:p>'

    $null = $r.InsertXML($newParaXml)

    Write-Host "Updated paragraph text:" $target.Range.Text
}
